{"js": "// 1) Split the run \"This is load data function\" into three runs:\n//    \"This is \" / \"a \" / \"load data function\"\n//    (done while the existing _GoBack bookmark is still in the paragraph so\n//    the inserted pieces are not re-coalesced into a single run).\nlet paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet loadDataPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"This is load data function\") {\n    loadDataPara = paras.items[i];\n    break;\n  }\n}\n\nconst wholeRange = loadDataPara.getRange(\"Whole\");\nwholeRange.insertText(\"\", \"Replace\");\nawait context.sync();\n\nwholeRange.insertText(\"This is \", \"Start\");\nawait context.sync();\nwholeRange.insertText(\"a \", \"End\");\nawait context.sync();\nwholeRange.insertText(\"load data function\", \"End\");\nawait context.sync();\n\n// 2) Remove the old \"_GoBack\" bookmark from that paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Delete the two now-empty list paragraphs that used to follow it.\nparas = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet afterIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"This is a load data function\") {\n    afterIdx = i;\n    break;\n  }\n}\n// The next two paragraphs are the empty ilvl=1 and ilvl=0 list items.\nparas.items[afterIdx + 1].delete();\nparas.items[afterIdx + 2].delete();\nawait context.sync();\n\n// 4) Remove the standalone empty paragraph right before \"Hyperparameters\".\nparas = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet hyperparamsIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Hyperparameters\") {\n    hyperparamsIdx = i;\n    break;\n  }\n}\nif (hyperparamsIdx > 0) {\n  paras.items[hyperparamsIdx - 1].delete();\n  await context.sync();\n}\n\n// 5) Split the \"LR \u2013 ...\" paragraph into two runs and insert a \"_GoBack\"\n//    bookmark between them, right after \"...loss with d\".\nparas = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet lrPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(\"LR\") === 0) {\n    lrPara = paras.items[i];\n    break;\n  }\n}\n\nconst found = lrPara.search(\"loss with d\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nconst splitPoint = found.items[0].getRange(\"After\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ParaIndexByText($doc, $exactText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs($i).Range.Text -eq ($exactText + \"`r\")) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Find-ParaIndexByPrefix($doc, $prefix) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs($i).Range.Text.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# 1) Split the run \"This is load data function\" into three runs:\n#    \"This is \" / \"a \" / \"load data function\".\n$idx = Find-ParaIndexByText $d \"This is load data function\"\n$target = $d.Paragraphs($idx)\n\n$rng = $target.Range\n$rng.End = $rng.End - 1\n$rng.Text = \"This is \"\n\n$insPt = $target.Range\n$insPt.End = $insPt.End - 1\n$insPt.Collapse(0)\n$insPt.InsertAfter(\"a \")\n\n$insPt2 = $target.Range\n$insPt2.End = $insPt2.End - 1\n$insPt2.Collapse(0)\n$insPt2.InsertAfter(\"load data function\")\n\n# 2) Remove the old \"_GoBack\" bookmark from that paragraph.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 3) Delete the two now-empty list paragraphs that used to follow it.\n$d.Paragraphs($idx + 1).Range.Delete()\n$d.Paragraphs($idx + 1).Range.Delete()\n\n# 4) Remove the standalone empty paragraph right before \"Hyperparameters\".\n$hpIdx = Find-ParaIndexByText $d \"Hyperparameters\"\n$d.Paragraphs($hpIdx - 1).Range.Delete()\n\n# 5) Split the \"LR \u2013 ...\" paragraph into two runs and insert a \"_GoBack\"\n#    bookmark between them, right after \"...loss with d\".\n$lrIdx = Find-ParaIndexByPrefix $d \"LR\"\n$lrPara = $d.Paragraphs($lrIdx)\n$splitRange = $lrPara.Range.Duplicate()\n$splitRange.Find.Execute(\"loss with d\")\n$splitRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $splitRange)\n"}
